$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = 9766.2900000000009
$ws.Range("B15").Value = 9820.2999999999993
$ws.Range("C15").Value = 78.48
$ws.Range("D15").Value = 78.05
$ws.Range("E15").Value = $false
$ws.Range("F15").Value = -0.55000000000000004
$ws.Range("G15").Value = 42624.611134259256
$ws.Range("G15").NumberFormat = "m/d/yy h:mm"
$ws.Range("H15").Value = $false
